# Update the "Countries" worksheet with new trial values and selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

# Brobdingnag (row 3): depth bound, max frontier size, and anytime version updated
$ws.Range("C3").Value = 90
$ws.Range("D3").Value = 600
$ws.Range("E3").Value = 1000

# Dinotopia (row 5): depth bound, max frontier size, and anytime version updated
$ws.Range("C5").Value = 85
$ws.Range("D5").Value = 550
$ws.Range("E5").Value = 1700

# Move the active selection to E5, matching the saved view state
$ws.Activate()
$ws.Range("E5").Select()
